$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I needed a bit more room for "All Members" entries.
$ws.Columns.Item(9).ColumnWidth = 10.6

# Row 19: Meeting 7 / All Members
$ws.Range("A19").Value = "Meeting 7"
$ws.Range("A19").Font.Bold = $true
$ws.Range("I19").Value = "All Members"
$ws.Range("I19").Font.Bold = $true

# Row 20: Problem sloving Session / All Members
$ws.Range("A20").Value = "Problem sloving Session"
$ws.Range("A20").Font.Bold = $true
$ws.Range("I20").Value = "All Members"
$ws.Range("I20").Font.Bold = $true

# Row 21: Work with functionality / 夏义
$ws.Range("A21").Value = "Work with functionality"
$ws.Range("A21").Font.Bold = $true
$ws.Range("I21").Value = "夏义"
$ws.Range("I21").Font.Bold = $true

$ws.Range("J21").Select() | Out-Null
